# Weekly refresh of the "Hortaliza, Terminal Hortofrutícola Agro Chillán - Ajo" sheet.
# Two new price records (week of 2023-10-12, serial 45211) are inserted at the top of the
# data block (row 438), pushing all existing data rows down by two. The table therefore
# grows from A1:R531 to A1:R533.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (above the former row 438),
# shifting the existing 438:531 block down to 440:533.
$ws.Range("A438:A439").EntireRow.Insert()

# --- New row 438: Ajo / Chino, $/caja 10 kilos ---
$ws.Cells.Item(438, 1).Value = 7
$ws.Cells.Item(438, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(438, 3).Value = "Ñuble"
$ws.Cells.Item(438, 4).Value = 45211
$ws.Cells.Item(438, 5).Value = 16
$ws.Cells.Item(438, 6).Value = 100112003
$ws.Cells.Item(438, 7).Value = "Ajo"
$ws.Cells.Item(438, 8).Value = "Chino"
$ws.Cells.Item(438, 9).Value = "Primera"
$ws.Cells.Item(438, 10).Value = 50
$ws.Cells.Item(438, 11).Value = 20000
$ws.Cells.Item(438, 12).Value = 20000
$ws.Cells.Item(438, 13).Value = 20000
$ws.Cells.Item(438, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(438, 15).Value = "China"
$ws.Cells.Item(438, 16).Value = 2000
$ws.Cells.Item(438, 17).Value = 10
$ws.Cells.Item(438, 18).Value = "Hortaliza"

# --- New row 439: Ajo / Chino, $/malla 10 kilos ---
$ws.Cells.Item(439, 1).Value = 7
$ws.Cells.Item(439, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(439, 3).Value = "Ñuble"
$ws.Cells.Item(439, 4).Value = 45211
$ws.Cells.Item(439, 5).Value = 16
$ws.Cells.Item(439, 6).Value = 100112003
$ws.Cells.Item(439, 7).Value = "Ajo"
$ws.Cells.Item(439, 8).Value = "Chino"
$ws.Cells.Item(439, 9).Value = "Primera"
$ws.Cells.Item(439, 10).Value = 50
$ws.Cells.Item(439, 11).Value = 22000
$ws.Cells.Item(439, 12).Value = 22000
$ws.Cells.Item(439, 13).Value = 22000
$ws.Cells.Item(439, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(439, 15).Value = "China"
$ws.Cells.Item(439, 16).Value = 2200
$ws.Cells.Item(439, 17).Value = 10
$ws.Cells.Item(439, 18).Value = "Hortaliza"
